$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows added below the existing data (rows 4-13 already populated)
$ws.Range("A14").Value = "AC Akım"
$ws.Range("A16").Value = "Şarj"
$ws.Range("A17").Value = "Float Şarj"
$ws.Range("B17").Value = "Devreye al. Akım voltajı kontrol et. "
$ws.Range("A20").Value = "Oto şarj test"

# Column width adjustments: A narrower, new B column added
# (values chosen land on the same quantized width the host stores for 13.22 / 29.44)
$ws.Columns.Item(1).ColumnWidth = 12.3
$ws.Columns.Item(2).ColumnWidth = 28.6

# Update selection to match final cursor position from the diff
$ws.Range("B21").Select()
